$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the "TreePoseIntro" row (old row 6) to hold the
# new "TreePoseOnlineID" / video-id pair. This pushes every row below it
# (and the hyperlink-carrying cells in column B) down by one.
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = "TreePoseOnlineID"
$ws.Range("B6").Value = "i3BPLAud0u0"

# Match the wrap-text formatting used by the other "value" cells in column B
# (e.g. B8, which holds the long TreePoseIntro description) without
# introducing a brand-new style entry.
$ws.Range("B8").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row heights settle slightly differently once the new row is present -
# mirror the values from the saved workbook.
$ws.Rows.Item(7).RowHeight = 87
$ws.Rows.Item(8).RowHeight = 275.5
$ws.Rows.Item(9).RowHeight = 43.5
$ws.Rows.Item(22).RowHeight = 362.5

# Row-insert doesn't shift the worksheet's stored hyperlink references, so
# rebuild the collection from scratch at the (now shifted-down-by-one)
# target cells. Stash the clean "Hyperlink" cell style first so it can be
# restored afterwards (re-adding a hyperlink otherwise forces a fresh style).
$ws.Range("B5").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B5"), "file:///\\Videos\TreePose.mp4")
$ws.Hyperlinks.Add($ws.Range("B11"), "file:///\\Videos\CatPose.mp4")
$ws.Hyperlinks.Add($ws.Range("B13"), "file:///\\Videos\BowPose.mp4")
$ws.Hyperlinks.Add($ws.Range("B15"), "file:///\\Videos\BridgePose.mp4")
$ws.Hyperlinks.Add($ws.Range("B17"), "file:///\\Videos\CamelPose.mp4")
$ws.Hyperlinks.Add($ws.Range("B19"), "file:///\\Videos\ChairPose.mp4")
$ws.Hyperlinks.Add($ws.Range("B21"), "file:///\\Videos\CowPose.mp4")
$ws.Hyperlinks.Add($ws.Range("B23"), "file:///\\Videos\TrianglePose.mp4")

$ws.Range("Z1").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Z1").Clear()

# Selection moves up to the row that now holds the long "steps" text.
[void]$ws.Range("B7").Select()
